# "Casos de prueba" workbook update - día 15
# - Rewrites the "Login requerido?" column (E) header text and values
#   (previously it held duplicate long text, now short SI/NO labels).
# - Reworks rows 8-22 ("Caso #8".."Caso #22") turning many entries into
#   generic reusable case numbers/actions, and rewrites several expected
#   results to be clearer ("Crear Producto, no se debe poder si el modelo
#   ya existe (solo admin)."  etc.) and fixes typos ("porducto" -> "producto").
# - Fills in 5 new test cases (#23-#27): view/edit profile, create/edit
#   avatar, and updated login/logout/register wording, plus a new
#   "Ver avatar" case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header + "Login requerido?" values for the already-existing
#     rows 6-13 (Casos #1-#7) ---
$ws.Range("E6").Value  = "Login requerido?"
$ws.Range("E7").Value  = "NO"
$ws.Range("E8").Value  = "NO"
$ws.Range("E9").Value  = "NO"
$ws.Range("E10").Value = "NO"
$ws.Range("E11").Value = "NO"
$ws.Range("E12").Value = "NO"
$ws.Range("E13").Value = "NO"

# --- Caso #8 (Registrar usuario) ---
$ws.Range("C14").Value = "Lograr registrar un usuario desde un forms propio, , ver página de registro correcto."
$ws.Range("D14").Value = 45275
$ws.Range("E14").Value = "NO"

# --- Caso #9 (Hacer un login) ---
$ws.Range("A15").Value = "Caso #9"
$ws.Range("B15").Value = "Hacer un login"
$ws.Range("C15").Value = "Lograr hacer login, ver página de login correcto."
$ws.Range("D15").Value = 45275
$ws.Range("E15").Value = "NO"

# --- Caso #10 (Hacer un logout) ---
$ws.Range("A16").Value = "Caso #10"
$ws.Range("B16").Value = "Hacer un logout"
$ws.Range("C16").Value = "Lograr hacer logout, ver página de logout correcto."
$ws.Range("D16").Value = 45275
$ws.Range("E16").Value = "NO"

# --- Caso #11 (Acceder a menú usuario) ---
$ws.Range("A17").Value = "Caso #11"
$ws.Range("B17").Value = "Acceder a menú usuario"
$ws.Range("C17").Value = "Ver el menu de usuario"
$ws.Range("E17").Value = "SI"

# --- Caso #12 (Acceder a menú admin) ---
$ws.Range("A18").Value = "Caso #12"
$ws.Range("B18").Value = "Acceder a menú admin"
$ws.Range("C18").Value = "Ver el menu de admin"
$ws.Range("E18").Value = "SI"

# --- Caso #13 (Crear Producto) ---
$ws.Range("A19").Value = "Caso #13"
$ws.Range("B19").Value = "Crear Producto"
$ws.Range("C19").Value = "Crear Producto, no se debe poder si el modelo ya existe (solo admin)."
$ws.Range("E19").Value = "SI"

# --- Caso #14 (Ver Productos) ---
$ws.Range("A20").Value = "Caso #14"
$ws.Range("B20").Value = "Ver Productos"
$ws.Range("C20").Value = "Ver lista de productos (admin)"
$ws.Range("E20").Value = "SI"

# --- Caso #15 (Editar producto) ---
$ws.Range("A21").Value = "Caso #15"
$ws.Range("B21").Value = "Editar producto"
$ws.Range("C21").Value = "Editar producto (solo admin)"
$ws.Range("E21").Value = "SI"

# --- Caso #16 (Borrar producto) ---
$ws.Range("A22").Value = "Caso #16"
$ws.Range("B22").Value = "Borrar producto"
$ws.Range("C22").Value = "Borrar producto (solo admin)"
$ws.Range("E22").Value = "SI"

# --- Caso #17 (Crear reseña) ---
$ws.Range("A23").Value = "Caso #17"
$ws.Range("B23").Value = "Crear reseña"
$ws.Range("C23").Value = "Crear reseña, no se debe poder si el modelo ya existe (solo admin)."
$ws.Range("E23").Value = "SI"

# --- Caso #18 (Ver reseña (admin)) ---
$ws.Range("A24").Value = "Caso #18"
$ws.Range("B24").Value = "Ver reseña (admin)"
$ws.Range("C24").Value = "Ver reseña desde el menu de admin"
$ws.Range("E24").Value = "SI"

# --- Caso #19 (Editar reseña) ---
$ws.Range("A25").Value = "Caso #19"
$ws.Range("B25").Value = "Editar reseña"
$ws.Range("C25").Value = "Editar reseña (solo admin)"
$ws.Range("E25").Value = "SI"

# --- Caso #20 (Eliminar reseña) ---
$ws.Range("A26").Value = "Caso #20"
$ws.Range("B26").Value = "Eliminar reseña"
$ws.Range("C26").Value = "Eliminar reseña (solo admin)"
$ws.Range("E26").Value = "SI"

# --- Caso #21 (Ingresar compra) ---
$ws.Range("A27").Value = "Caso #21"
$ws.Range("B27").Value = "Ingresar compra"
$ws.Range("C27").Value = "Ingrersar compra desde menu de usuario"
$ws.Range("E27").Value = "SI"

# --- Caso #22 (Ver compras usuario) ---
$ws.Range("A28").Value = "Caso #22"
$ws.Range("B28").Value = "Ver compras usuario"
$ws.Range("C28").Value = "Ver compras hechas por el usuario logueado"
$ws.Range("E28").Value = "SI"

# --- Caso #23 (Ver compras admin) - was a blank template row ---
$ws.Range("A29").Value = "Caso #23"
$ws.Range("B29").Value = "Ver compras admin"
$ws.Range("C29").Value = "Ver las compras de todos los usuarios"
$ws.Range("D29").Value = 45274
$ws.Range("E29").Value = "SI"
$ws.Range("F29").Value = "-"
$ws.Range("G29").Value = "OK"

# --- Caso #24 (Ver Perfil) - new ---
$ws.Range("A30").Value = "Caso #24"
$ws.Range("B30").Value = "Ver Perfil"
$ws.Range("C30").Value = "Ver perfil activo"
$ws.Range("D30").Value = 45275
$ws.Range("E30").Value = "SI"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "OK"

# --- Caso #25 (Editar Perfil) - new ---
$ws.Range("A31").Value = "Caso #25"
$ws.Range("B31").Value = "Editar Perfil"
$ws.Range("C31").Value = "Editar perfil activo"
$ws.Range("D31").Value = 45275
$ws.Range("E31").Value = "SI"
$ws.Range("F31").Value = "-"
$ws.Range("G31").Value = "OK"

# --- Caso #26 (Crear o editar avatar) - new ---
$ws.Range("A32").Value = "Caso #26"
$ws.Range("B32").Value = "Crear o editar avatar"
$ws.Range("C32").Value = "Crear o editar avatar desde menu de usuario o admin"
$ws.Range("D32").Value = 45275
$ws.Range("E32").Value = "SI"
$ws.Range("F32").Value = "-"
$ws.Range("G32").Value = "OK"

# --- Caso #27 (Ver avatar) - new ---
$ws.Range("A33").Value = "Caso #27"
$ws.Range("B33").Value = "Ver avatar"
$ws.Range("C33").Value = "Se debe ver el avatar en todas las p{aginas del proyecto"
$ws.Range("D33").Value = 45275
$ws.Range("E33").Value = "SI"
$ws.Range("F33").Value = "Hasta el momento se ven, falta hacer la parte de blog"
$ws.Range("G33").Value = "OK"

# Leave the selection where the author ended up editing.
$ws.Range("G33").Select()
